$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Constants used with PasteSpecial ---
$xlPasteValues  = -4163
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) "Ativação:" date changes from 01/01/2012 to 01/01/2023.
#    Four cells (B8, C8, B13, C13) currently share that string.
#    Assigning a literal "01/01/2023" via .Value would be auto-coerced into a
#    serial date by the engine's type inference, changing both the cell type
#    and its style. To keep the cells as plain shared-string text with their
#    existing style untouched, we push the text through a text formula and
#    then flatten it back to a static value with Copy/PasteSpecial(values).
# ---------------------------------------------------------------------------
foreach ($addr in @("B8", "C8", "B13", "C13")) {
    $rng = $ws.Range($addr)
    $rng.Formula = "=""01/01/2023"""
    $rng.Copy()
    $rng.PasteSpecial($xlPasteValues)
}

# ---------------------------------------------------------------------------
# 2) New seminar-related sentences added in three places. Each new line adds
#    matching red/black columns (B = normal, C = "modified" red copy), mirror
#    -ing the existing B/C pattern elsewhere in the sheet. We seed the text
#    via .Value (safe here - none of this text is a date), then copy the
#    number/font formatting across from the nearest existing B/C pair so the
#    new cells pick up style s="2" (B) / s="3" (C) without creating any new
#    style entries.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row = 11; Text = "Provide student seminars on current topics in Physics, Technology and Engineering." },
    @{ Row = 14; Text = "Seminars covering the current and future scenarios of the high technology industry and the field of activity of the physical engineer." },
    @{ Row = 16; Text = "Seminars followed by debates with professionals and undergraduate and graduate students on relevant and current topics in the areas of Physics, Technology and Engineering, ranging from basic research to the industrial and services segment." }
)

foreach ($item in $newRows) {
    $row = $item.Row

    $bCell = $ws.Range("B$row")
    $bCell.Value = $item.Text
    $ws.Range("B10").Copy()
    $bCell.PasteSpecial($xlPasteFormats)

    $cCell = $ws.Range("C$row")
    $cCell.Value = $item.Text
    $ws.Range("C10").Copy()
    $cCell.PasteSpecial($xlPasteFormats)
}
